# 自动更新Excel文件 - 2025-11-29 23:12:09
# Decrement the "剩余" (remaining) value in column E by 1 for every data
# row (rows 2-99), except row 36 whose value stays unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 99; $row++) {
    if ($row -eq 36) {
        continue
    }
    $cell = $ws.Cells.Item($row, 5)
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value2 = $current - 1
    }
}
